$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.2781022865678179
$ws.Range("D2").Value = 0.7826181163526409

# Row 3
$ws.Range("C3").Value = -0.7297671176796202
$ws.Range("D3").Value = 0.4705308036228892

# Row 4
$ws.Range("C4").Value = -1.699897949074796
$ws.Range("D4").Value = 0.09828244218043558
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = -2.099816983818729
$ws.Range("D5").Value = 0.04324053305676823

# Row 6
$ws.Range("C6").Value = -0.2771445933377392
$ws.Range("D6").Value = 0.7833471780010841

# Row 7
$ws.Range("C7").Value = -1.255393721793903
$ws.Range("D7").Value = 0.2178956497620161

# Row 8
$ws.Range("C8").Value = -1.760255249330729
$ws.Range("D8").Value = 0.08735908148987104
$ws.Range("G8").Value = "No"

# Row 9
$ws.Range("C9").Value = -1.518638694042757
$ws.Range("D9").Value = 0.1380984426700713

# Row 10
$ws.Range("C10").Value = -1.938158091634261
$ws.Range("D10").Value = 0.06094376260187717
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = -0.3568982988580892
$ws.Range("D11").Value = 0.7233735659337914
